$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check on" timestamp in F1
$ws.Range("F1").Value = "Last status check on: 03.02.2022 11:30"

# D3: change from text "+1.0" to numeric value 1
$ws.Range("D3").Value = 1

# E3: change from text timestamp to a real date/time value, formatted like the
# other date cells in column E (style index 2 / numFmt "YYYY-MM-DD HH:MM:SS")
$ws.Range("E3").NumberFormat = $ws.Range("E2").NumberFormat
$ws.Range("E3").Value = 44595.4703125
